$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixes: replace stray commas with periods in Razon Social / Nombre Fantasia ---
# (Plain .Value assignment: text already fails numeric parsing, so Excel keeps it as a
# shared string with no number-format side effects.)
$ws.Range("E16").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F16").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E32").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E33").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F33").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E34").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E48").Value = 'RICCOTTI. MARIANA EDITH'
$ws.Range("E58").Value = 'GIMENEZ. ROBERTO ADRIAN'
$ws.Range("F58").Value = 'GIMENEZ. ROBERTO ADRIAN'
$ws.Range("E65").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'

# --- Amount column (Importe): the scraped text used es-AR grouping - "1.234,56" - and needs to
# become plain "1234.56" text (same string type, just reformatted). Re-typing the digits via
# .Value would make Excel auto-convert the now-valid-looking number into a real numeric cell,
# so the column is first marked as Text ("@") and the substitution is performed in place with
# Replace (remove the thousands "." separator, then turn the decimal "," into "."); the format
# is switched back to the original Normal style afterwards so no visible formatting changes.
$amountRange = $ws.Range("H2:H86")
$amountRange.NumberFormat = "@"
[void]$amountRange.Replace(".", "")
[void]$amountRange.Replace(",", ".")
$amountRange.Style = "Normal"
